$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (company_name "6") ---
$ws.Range("D2").Value = 0.067
$ws.Range("E2").Value = -0.0471
$ws.Range("G2").Value = 0.08506218135035974
$ws.Range("H2").Value = 0.08506218135035974
$ws.Range("I2").Value = 0.05796445418592461
$ws.Range("J2").Value = 0.04920230452675968
$ws.Range("K2").Value = 12.948
$ws.Range("L2").Value = 0.0616953352075094
$ws.Range("M2").Value = 5.348
$ws.Range("N2").Value = 0.04788253200823708
$ws.Range("O2").Value = 0.4130367624343528
$ws.Range("P2").Value = 5.348
$ws.Range("Q2").Value = 0.04788253200823708
$ws.Range("R2").Value = 0.4130367624343528
$ws.Range("U2").Value = 10.973
$ws.Range("V2").Value = 0.09824514280598085
$ws.Range("W2").Value = 0.05364028213166144
$ws.Range("X2").Value = 0.0522719173994006
$ws.Range("Y2").Value = 0.001368364732260846
$ws.Range("Z2").Value = 1.382515497058688
$ws.Range("AA2").Value = 0.04305590897374717
$ws.Range("AB2").Value = 0.0522719173994006
$ws.Range("AC2").Value = -0.00921600842565343
$ws.Range("AD2").Value = 1.302
$ws.Range("AF2").Value = 1.302
$ws.Range("AG2").Value = -9.671000000000001
$ws.Range("AH2").Value = 0.01152293967714529
$ws.Range("AI2").Value = 0.008443470253304108
$ws.Range("AJ2").Value = -0.09479606739920997
$ws.Range("AK2").Value = -0.06752124220653639
$ws.Range("AL2").Value = 0.171
$ws.Range("AM2").Value = 0.171
$ws.Range("AN2").Value = 0.08651162790697675
$ws.Range("AO2").Value = 71.14035087719299
$ws.Range("AP2").Value = -0.642591362126246
$ws.Range("AQ2").Value = 71.14035087719299

# --- Row 3: renamed to PT Lippo General Insurance Tbk (IDX:LPGI) ---
$ws.Range("B3").Value = "PT Lippo General Insurance Tbk (IDX:LPGI)"
$ws.Range("D3").Value = 0.0603
$ws.Range("E3").Value = -0.0171
$ws.Range("G3").Value = 0.08562577447335812
$ws.Range("H3").Value = 0.08562577447335812
$ws.Range("I3").Value = 0.09752168525402725
$ws.Range("J3").Value = 0.08635387936203381
$ws.Range("K3").Value = 7.14
$ws.Range("L3").Value = 0.08847583643122676
$ws.Range("M3").Value = 2.16
$ws.Range("N3").Value = 0.05901639344262295
$ws.Range("O3").Value = 0.3025210084033614
$ws.Range("P3").Value = 2.16
$ws.Range("Q3").Value = 0.05901639344262295
$ws.Range("R3").Value = 0.3025210084033614
$ws.Range("V3").Value = 0.04562841530054645
$ws.Range("W3").Value = 0.1216354344122657
$ws.Range("X3").Value = 0.0522719173994006
$ws.Range("Y3").Value = 0.06936351701286515
$ws.Range("Z3").Value = 1.416037901386208
$ws.Range("AA3").Value = 0.1222803661083722
$ws.Range("AB3").Value = 0.0522719173994006
$ws.Range("AC3").Value = 0.07000844870897156
$ws.Range("AJ3").Value = -0.04780990552533639
$ws.Range("AK3").Value = -0.03266184236260512
$ws.Range("AL3").Value = 0.007
$ws.Range("AM3").Value = 0.007
$ws.Range("AO3").Value = 1124.285714285714
$ws.Range("AP3").Value = -0.1976331360946746
$ws.Range("AQ3").Value = 1124.285714285714

# --- Row 4: PT Asuransi Ramayana Tbk (IDX:ASRM) ---
$ws.Range("D4").Value = 0.14
$ws.Range("E4").Value = -0.0471
$ws.Range("G4").Value = 0.05968289920724801
$ws.Range("H4").Value = 0.05968289920724801
$ws.Range("I4").Value = 0.0304643261608154
$ws.Range("J4").Value = 0.02645620189262338
$ws.Range("K4").Value = 4.1
$ws.Range("L4").Value = 0.0464326160815402
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 0.02989130434782609
$ws.Range("O4").Value = 0.2682926829268293
$ws.Range("P4").Value = 1.1
$ws.Range("Q4").Value = 0.02989130434782609
$ws.Range("R4").Value = 0.2682926829268293
$ws.Range("U4").Value = 4.69
$ws.Range("V4").Value = 0.1274456521739131
$ws.Range("W4").Value = 0.1339869281045752
$ws.Range("X4").Value = 0.05350268916404648
$ws.Range("Y4").Value = 0.08048423894052868
$ws.Range("Z4").Value = 3.175116864437253
$ws.Range("AA4").Value = 0.08400153279822524
$ws.Range("AB4").Value = 0.05275673456360992
$ws.Range("AC4").Value = 0.03124479823461532
$ws.Range("AD4").Value = 1.24
$ws.Range("AF4").Value = 1.24
$ws.Range("AG4").Value = -3.45
$ws.Range("AH4").Value = 0.03259726603575184
$ws.Range("AI4").Value = 0.03719256148770245
$ws.Range("AJ4").Value = -0.103448275862069
$ws.Range("AK4").Value = -0.1204188481675393
$ws.Range("AL4").Value = 0.154
$ws.Range("AM4").Value = 0.154
$ws.Range("AN4").Value = 0.3280423280423281
$ws.Range("AO4").Value = 17.46753246753247
$ws.Range("AP4").Value = -0.9126984126984128
$ws.Range("AQ4").Value = 17.46753246753247

# --- Row 5: PT Victoria Insurance Tbk (IDX:VINS) ---
$ws.Range("D5").Value = 0.09949999999999999
$ws.Range("E5").Value = 0.0102
$ws.Range("G5").Value = 0.599250936329588
$ws.Range("H5").Value = 0.599250936329588
$ws.Range("I5").Value = 0.2164794007490637
$ws.Range("J5").Value = 0.2120454853120347
$ws.Range("K5").Value = 0.8129999999999999
$ws.Range("L5").Value = 0.3044943820224719
$ws.Range("M5").Value = 0.9330000000000001
$ws.Range("N5").Value = 0.09749216300940439
$ws.Range("O5").Value = 1.14760147601476
$ws.Range("P5").Value = 0.9330000000000001
$ws.Range("Q5").Value = 0.09749216300940439
$ws.Range("R5").Value = 1.14760147601476
$ws.Range("U5").Value = 0.229
$ws.Range("V5").Value = 0.02392894461859979
$ws.Range("W5").Value = 0.06159090909090909
$ws.Range("X5").Value = 0.0522719173994006
$ws.Range("Y5").Value = 0.009318991691508495
$ws.Range("Z5").Value = 0.2036613272311213
$ws.Range("AA5").Value = 0.04318546497201622
$ws.Range("AB5").Value = 0.0522719173994006
$ws.Range("AC5").Value = -0.00908645242738438
$ws.Range("AG5").Value = -0.229
$ws.Range("AJ5").Value = -0.02451557649073975
$ws.Range("AK5").Value = -0.01912956311085122
$ws.Range("AP5").Value = -0.3669871794871795

# --- Row 6: renamed to PT Asuransi Dayin Mitra Tbk (IDX:ASDM) ---
$ws.Range("B6").Value = "PT Asuransi Dayin Mitra Tbk (IDX:ASDM)"
$ws.Range("D6").Value = -0.0184
$ws.Range("E6").Value = -0.189
$ws.Range("G6").Value = 0.2395161290322581
$ws.Range("H6").Value = 0.2395161290322581
$ws.Range("I6").Value = 0.09596774193548387
$ws.Range("J6").Value = 0.07453261125500353
$ws.Range("K6").Value = 1.06
$ws.Range("L6").Value = 0.08548387096774193
$ws.Range("M6").Value = 1.01
$ws.Range("N6").Value = 0.08211382113821138
$ws.Range("O6").Value = 0.9528301886792453
$ws.Range("P6").Value = 1.01
$ws.Range("Q6").Value = 0.08211382113821138
$ws.Range("R6").Value = 0.9528301886792453
$ws.Range("U6").Value = 1.16
$ws.Range("V6").Value = 0.09430894308943089
$ws.Range("W6").Value = 0.04568965517241379
$ws.Range("X6").Value = 0.0522719173994006
$ws.Range("Y6").Value = -0.006582262226986803
$ws.Range("Z6").Value = 0.575940548072457
$ws.Range("AA6").Value = 0.04292635297547812
$ws.Range("AB6").Value = 0.0522719173994006
$ws.Range("AC6").Value = -0.009345564423922481
$ws.Range("AG6").Value = -1.16
$ws.Range("AJ6").Value = -0.104129263913824
$ws.Range("AK6").Value = -0.05566218809980806
$ws.Range("AP6").Value = -0.7341772151898733

# --- Row 7: renamed to PT Asuransi Bintang Tbk (IDX:ASBI) ---
$ws.Range("B7").Value = "PT Asuransi Bintang Tbk (IDX:ASBI)"
$ws.Range("D7").Value = 0.0737
$ws.Range("E7").Value = -0.27
$ws.Range("G7").Value = 0.04508771929824561
$ws.Range("H7").Value = 0.04508771929824561
$ws.Range("I7").Value = 0.01684210526315789
$ws.Range("J7").Value = 0.009817909703167871
$ws.Range("K7").Value = 0.246
$ws.Range("L7").Value = 0.0143859649122807
$ws.Range("M7").Value = 0.145
$ws.Range("N7").Value = 0.01863753213367609
$ws.Range("O7").Value = 0.589430894308943
$ws.Range("P7").Value = 0.145
$ws.Range("Q7").Value = 0.01863753213367609
$ws.Range("R7").Value = 0.589430894308943
$ws.Range("U7").Value = 2.38
$ws.Range("V7").Value = 0.3059125964010282
$ws.Range("W7").Value = 0.01223880597014925
$ws.Range("X7").Value = 0.0525629996676459
$ws.Range("Y7").Value = -0.04032419369749665
$ws.Range("Z7").Value = 0.919206579583938
$ws.Range("AA7").Value = 0.009024687196912896
$ws.Range("AB7").Value = 0.0524619490553973
$ws.Range("AC7").Value = -0.0434372618584844
$ws.Range("AD7").Value = 0.062
$ws.Range("AF7").Value = 0.062
$ws.Range("AG7").Value = -2.318
$ws.Range("AH7").Value = 0.007906146391226727
$ws.Range("AI7").Value = 0.003153290611331502
$ws.Range("AJ7").Value = -0.4243866715488832
$ws.Range("AK7").Value = -0.1341279944450874
$ws.Range("AL7").Value = 0.01
$ws.Range("AM7").Value = 0.01
$ws.Range("AN7").Value = 0.07989690721649484
$ws.Range("AO7").Value = 28.8
$ws.Range("AP7").Value = -2.987113402061856
$ws.Range("AQ7").Value = 28.8

# --- Row 8: renamed to PT Asuransi Jasa Tania Tbk (IDX:ASJT) ---
$ws.Range("B8").Value = "PT Asuransi Jasa Tania Tbk (IDX:ASJT)"
$ws.Range("D8").Value = -0.09
$ws.Range("G8").Value = 0.03804597701149426
$ws.Range("H8").Value = 0.03804597701149426
$ws.Range("I8").Value = -0.05183908045977012
$ws.Range("J8").Value = -0.05183908045977012
$ws.Range("K8").Value = -0.411
$ws.Range("L8").Value = -0.04724137931034483
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 0
$ws.Range("U8").Value = 0.844
$ws.Range("V8").Value = 0.09768518518518518
$ws.Range("W8").Value = -0.02777027027027027
$ws.Range("X8").Value = 0.0522719173994006
$ws.Range("Y8").Value = -0.08004218766967086
$ws.Range("Z8").Value = 0.632267441860465
$ws.Range("AA8").Value = -0.03277616279069767
$ws.Range("AB8").Value = 0.0522719173994006
$ws.Range("AC8").Value = -0.08504808019009827
$ws.Range("AD8").Value = 0
$ws.Range("AF8").Value = 0
$ws.Range("AG8").Value = -0.844
$ws.Range("AH8").Value = 0
$ws.Range("AI8").Value = 0
$ws.Range("AJ8").Value = -0.1082606464853771
$ws.Range("AK8").Value = -0.0631925726265349
$ws.Range("AL8").Value = 0
$ws.Range("AM8").Value = 0
$ws.Range("AN8").Value = 0
$ws.Range("AP8").Value = 5.274999999999999

# --- Clear cells removed in row 8 (no longer applicable) ---
$ws.Range("E8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("AO8").ClearContents()
$ws.Range("AQ8").ClearContents()

Write-Host "Done updating capital structure database"
